# Generate Report for Handback
# Updates the localization-status report after a handback: marks both
# target languages as "Handed back: in sync with en-US", records the
# generated handback xliff file name + timestamp for zh-cn and de-de,
# links the target markdown file, and widens the columns that now need
# to show the longer status text / file names.

$wb = $excel.ActiveWorkbook

$wsOverview = $wb.Worksheets.Item("Overview")
$wsZhCn     = $wb.Worksheets.Item("zh-cn")
$wsDeDe     = $wb.Worksheets.Item("de-de")

$statusText = "Handed back: in sync with en-US"

# ---------------------------------------------------------------------
# Overview sheet: status text for each language column
# ---------------------------------------------------------------------
$wsOverview.Range("E2").Value = $statusText
$wsOverview.Range("F2").Value = $statusText
$wsOverview.Range("E3").Value = $statusText
$wsOverview.Range("F3").Value = $statusText

# Widen the zh-cn / de-de status columns so the longer text is visible.
$wsOverview.Range("E1").ColumnWidth = 29.17
$wsOverview.Range("F1").ColumnWidth = 29.17

# ---------------------------------------------------------------------
# zh-cn sheet
# ---------------------------------------------------------------------
$wsZhCn.Range("C2").Value = $statusText
$wsZhCn.Range("C3").Value = $statusText

$wsZhCn.Hyperlinks.Add($wsZhCn.Cells.Item(2, 9), "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/c6b94f2e15aca95e9f021b41c6043dff12c2fce4/e2e/a.md", "", "", "a.md")
$wsZhCn.Cells.Item(2, 9).Font.Name = "Calibri"
$wsZhCn.Cells.Item(2, 9).Font.Size = 11
$wsZhCn.Cells.Item(2, 9).Font.Underline = 2
$wsZhCn.Cells.Item(2, 9).Font.Color = 15570276

$wsZhCn.Hyperlinks.Add($wsZhCn.Cells.Item(3, 9), "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/c6b94f2e15aca95e9f021b41c6043dff12c2fce4/e2e/a.md", "", "", "a.md")
$wsZhCn.Cells.Item(3, 9).Font.Name = "Calibri"
$wsZhCn.Cells.Item(3, 9).Font.Size = 11
$wsZhCn.Cells.Item(3, 9).Font.Underline = 2
$wsZhCn.Cells.Item(3, 9).Font.Color = 15570276

$wsZhCn.Range("J2").Value = "a.6631f68b315a3f7ddcdc141802fdb6ef151d1df2.zh-cn.xlf"
$wsZhCn.Range("J3").Value = "a.6631f68b315a3f7ddcdc141802fdb6ef151d1df2.zh-cn.xlf"

$wsZhCn.Range("K2").Value = "2016-09-05 06:42:35"
$wsZhCn.Range("K3").Value = "2016-09-05 06:42:35"

$wsZhCn.Range("C1").ColumnWidth = 29.17
$wsZhCn.Range("J1").ColumnWidth = 39.17

# ---------------------------------------------------------------------
# de-de sheet
# ---------------------------------------------------------------------
$wsDeDe.Range("C2").Value = $statusText
$wsDeDe.Range("C3").Value = $statusText

$wsDeDe.Hyperlinks.Add($wsDeDe.Cells.Item(2, 9), "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/c6b94f2e15aca95e9f021b41c6043dff12c2fce4/e2e/a.md", "", "", "a.md")
$wsDeDe.Cells.Item(2, 9).Font.Name = "Calibri"
$wsDeDe.Cells.Item(2, 9).Font.Size = 11
$wsDeDe.Cells.Item(2, 9).Font.Underline = 2
$wsDeDe.Cells.Item(2, 9).Font.Color = 15570276

$wsDeDe.Hyperlinks.Add($wsDeDe.Cells.Item(3, 9), "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/c6b94f2e15aca95e9f021b41c6043dff12c2fce4/e2e/a.md", "", "", "a.md")
$wsDeDe.Cells.Item(3, 9).Font.Name = "Calibri"
$wsDeDe.Cells.Item(3, 9).Font.Size = 11
$wsDeDe.Cells.Item(3, 9).Font.Underline = 2
$wsDeDe.Cells.Item(3, 9).Font.Color = 15570276

$wsDeDe.Range("J2").Value = "a.6631f68b315a3f7ddcdc141802fdb6ef151d1df2.de-de.xlf"
$wsDeDe.Range("J3").Value = "a.6631f68b315a3f7ddcdc141802fdb6ef151d1df2.de-de.xlf"

$wsDeDe.Range("K2").Value = "2016-09-05 06:42:42"
$wsDeDe.Range("K3").Value = "2016-09-05 06:42:42"

$wsDeDe.Range("C1").ColumnWidth = 29.17
$wsDeDe.Range("J1").ColumnWidth = 39.17
